# Simplify and clarify scenario names
# Rewrites the verbose LOCATION descriptions (column C) in the link table
# down to short place names, resizes column C to fit, and leaves the
# cursor parked on C11 (matching the author's final selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "Hite"
$ws.Range("C3").Value  = "King Top"
$ws.Range("C4").Value  = "Bountiful"
$ws.Range("C5").Value  = "West of Green River"
$ws.Range("C6").Value  = "Orem"
$ws.Range("C7").Value  = "Rush Valley"
$ws.Range("C8").Value  = "Beaver Canyon"
$ws.Range("C9").Value  = "Snow Canyon"
$ws.Range("C10").Value = "Utah / Juab county line"
$ws.Range("C11").Value = "North of Cove Fort"
$ws.Range("C12").Value = "North of Zion"
$ws.Range("C13").Value = "East of Strawberry Reservoir"
$ws.Range("C14").Value = "Price Canyon"
$ws.Range("C15").Value = "West Bountiful"
$ws.Range("C16").Value = "Francis"
$ws.Range("C17").Value = "American Fork Canyon"
$ws.Range("C18").Value = "Cedar Canyon"
$ws.Range("C19").Value = "Weber Canyon"
$ws.Range("C20").Value = "Emigration Canyon"
$ws.Range("C21").Value = "Hyrum"
$ws.Range("C22").Value = "Box Elder Canyon"
$ws.Range("C23").Value = "Kingston"
$ws.Range("C24").Value = "Logan Canyon"
$ws.Range("C25").Value = "in Capitol Reef National Park"
$ws.Range("C26").Value = "near Bluffdale"
$ws.Range("C27").Value = "between Helper & Duchesne"
$ws.Range("C28").Value = "West of Hanksville"
$ws.Range("C29").Value = "Parley's Canyon"
$ws.Range("C30").Value = "Utah / Salt Lake county line"
$ws.Range("C31").Value = "SLC 1300 E"
$ws.Range("C32").Value = "SLC  2100 S"
$ws.Range("C33").Value = "Taylorsville"
$ws.Range("C34").Value = "West Valley City"
$ws.Range("C35").Value = "Cottonwood Heights"
$ws.Range("C36").Value = "West Jordan"
$ws.Range("C37").Value = "Arizona state line"
$ws.Range("C38").Value = "Provo Canyon"
$ws.Range("C39").Value = "Spanish Fork Canyon"
$ws.Range("C40").Value = "Colorado state line"
$ws.Range("C41").Value = "East of Cove Fort"
$ws.Range("C42").Value = "Salt Lake / Tooele county line"

# Widen column C (LOCATION) to fit the longest entry; ColumnWidth uses
# Excel's "characters" unit, which renders ~0.8333 wider in the saved
# XML width attribute, so back that padding out to land on 43.5.
$ws.Columns.Item(3).ColumnWidth = 42.666666666666664

# Park the active cell/selection on C11, matching the author's last
# on-screen selection when the workbook was saved.
$ws.Range("C11").Select()
